# Architecture diagram cleanup: straighten the connector lines and
# nudge a couple of shapes so they line up, plus the incidental
# "date last saved" placeholder bump that PowerPoint stamps on the
# slide master / layouts (9/7/2017 -> 9/12/2017).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Update the auto date placeholders on the slide master and on
#    every slide layout (type 16 = ppPlaceholderDate).
# ---------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16 -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "9/7/2017") {
                $sh.TextFrame.TextRange.Text = "9/12/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------
# 2) Straighten out the architecture-diagram shapes on slide 1.
# ---------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if ($sh.Name -eq "Rounded Rectangle 14") {
        $sh.Left = 231.25818897637797
    }
    elseif ($sh.Name -eq "Picture 18") {
        $sh.Left = 511.8955905511811
        $sh.Top = 36.34590751181103
    }
    elseif ($sh.Name -eq "Straight Connector 91") {
        $sh.Left = 397.6562992125984
        $sh.Top = 83.868031496063
        $sh.Width = 124.42716635433071
        $sh.Height = 68.97606299212599
    }
    elseif ($sh.Name -eq "Straight Connector 94") {
        $sh.Left = 534.3299262598425
        $sh.Top = 85.50299512598426
        $sh.Width = 0.718503937007874
        $sh.Height = 128.69268116535432
    }
}
